$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 77

# The date-like text would otherwise be auto-converted to a real date by
# Excel's type inference, so force it in as literal text (leading
# apostrophe), then drop the resulting "quote prefix" cell format so the
# cell ends up with the same (default) style as its neighbours.
$ws.Cells.Item($row, 1).Value = "'02/09/2026"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = 9365.299999999999
$ws.Cells.Item($row, 3).Value = 0.2422875806110708
$ws.Cells.Item($row, 4).Value = 0.7577124193889292
$ws.Cells.Item($row, 5).Value = -316.94
$ws.Cells.Item($row, 6).Value = -37.96
$ws.Cells.Item($row, 7).Value = -23795.55
$ws.Cells.Item($row, 8).Value = -77.03
$ws.Cells.Item($row, 9).Value = -1083.75
$ws.Cells.Item($row, 10).Value = -32.32
$ws.Cells.Item($row, 11).Value = -24879.3
$ws.Cells.Item($row, 12).Value = -72.65000000000001
